$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column G (rows 3-8): new numeric timings ----
$ws.Range("G3").Value = 880.06517899999994
$ws.Range("G4").Value = 1220.598821
$ws.Range("G5").Value = 11897.80673
$ws.Range("G6").Value = 14903.132895000001
$ws.Range("G7").Value = 15031.086153
$ws.Range("G8").Value = 20598.70709

# Rows 7:8 in column G need to pick up the right-aligned integer style
# already used by rows 3:6 (style index 9 in the original workbook).
$ws.Range("G7:G8").NumberFormat = "0"
$ws.Range("G7:G8").HorizontalAlignment = -4152

# ---- Column G/H/I (rows 9-13): now unavailable -> "N/A" ----
$ws.Range("G9:I13").Value = "N/A"
$ws.Range("G9:I13").NumberFormat = "0"
$ws.Range("G9:I13").HorizontalAlignment = -4152

# Rows 14-17 already show "N/A" but their G/H/I formatting is unified too.
$ws.Range("G14:I17").NumberFormat = "0"
$ws.Range("G14:I17").HorizontalAlignment = -4152

# ---- Column J (rows 3-17): new numeric timings ----
$ws.Range("J3").Value = 668.51707699999997
$ws.Range("J4").Value = 641.07298300000002
$ws.Range("J5").Value = 5497.9921039999999
$ws.Range("J6").Value = 6627.10772
$ws.Range("J7").Value = 7389.4488970000002
$ws.Range("J8").Value = 9962.8605580000003
$ws.Range("J9").Value = 17500.605729999999
$ws.Range("J10").Value = 26263.275151000002
$ws.Range("J11").Value = 36265.528225000002
$ws.Range("J12").Value = 41516.904816000002
$ws.Range("J13").Value = 54074.672651000001
$ws.Range("J14").Value = 67932.044580999995
$ws.Range("J15").Value = 98069.277424999993
$ws.Range("J16").Value = 124797.856075
$ws.Range("J17").Value = 138404.599919

# ---- Selection / view state ----
$ws.Range("J3:J17").Select() | Out-Null
